# Generate Report for Handoff
#
# Updates the "Latest Handoff Datetime" values for the row corresponding to
# file "02f42fc4-4350-425b-bb84-70aa4b7534e3.md" across the Overview,
# zh-cn and de-de sheets, reflecting a fresh handoff report generation.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Row 5 on each sheet corresponds to 02f42fc4-4350-425b-bb84-70aa4b7534e3.md
$zhcn.Range("E5").Value = "2016-03-21 12:33:43"
$dede.Range("E5").Value = "2016-03-21 12:33:47"
$overview.Range("D5").Value = "2016-03-21 12:33:47"
